$wb = $excel.ActiveWorkbook

# --- Sheet 1: question_answers ---
$ws1 = $wb.Worksheets.Item("question_answers")

$answers = @{
    5  = "3"
    8  = "2"
    9  = "4"
    10 = "5"
    11 = "2"
    12 = "2"
    13 = "3"
    14 = "4"
    15 = "5"
    16 = "4"
    17 = "4"
    18 = "2"
    19 = "3"
    20 = "5"
    23 = "5"
    24 = "5"
    25 = "6"
    26 = "6"
    27 = "6"
    28 = "1"
    29 = "5"
    30 = "4"
    31 = "3"
    32 = "4"
    33 = "3"
    34 = "2"
    35 = "1"
    36 = "3"
    37 = "6"
}

foreach ($row in $answers.Keys) {
    $cell = $ws1.Range("B$row")
    $cell.NumberFormat = "@"
    $cell.Value = $answers[$row]
}

# --- Sheet 2: outputs ---
$ws2 = $wb.Worksheets.Item("outputs")

# B2 is a genuine number
$ws2.Range("B2").Value = 15

# B4, B5, B6, B8, B9 are text values (even though they look numeric)
$textVals = @{
    4 = "5"
    5 = "10"
    6 = "25"
    8 = "74"
    9 = "87"
}
foreach ($row in $textVals.Keys) {
    $cell = $ws2.Range("B$row")
    $cell.NumberFormat = "@"
    $cell.Value = $textVals[$row]
}
